$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2535.75
$ws.Range("I125").Value = 3804.7144
$ws.Range("K125").Value = 34242.4296
$ws.Range("M125").Value = -31782.4296
$ws.Range("H137").Value = 1346.1714
$ws.Range("I137").Value = 1490.3889
$ws.Range("J137").Value = 1193.4706
$ws.Range("K137").Value = 4471.1667
$ws.Range("L137").Value = 3580.4118
$ws.Range("M137").Value = -1921.1667
$ws.Range("N137").Value = -8680.4118
$ws.Range("H138").Value = 3051.3132
$ws.Range("I138").Value = 1521.9117
$ws.Range("J138").Value = 4112.531
$ws.Range("K138").Value = 4565.7351
$ws.Range("L138").Value = 12337.593
$ws.Range("M138").Value = 574.2649000000001
$ws.Range("N138").Value = -22617.593

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14276.44
$ws.Range("I32").Value = 10765.731
$ws.Range("J32").Value = 21404.242
$ws.Range("K32").Value = 10765.731
$ws.Range("L32").Value = 21404.242
$ws.Range("M32").Value = -10478.731
$ws.Range("N32").Value = -21978.242
$ws.Range("H45").Value = 151206.86
$ws.Range("I45").Value = 263124.38
$ws.Range("J45").Value = 1983.5
$ws.Range("K45").Value = 263124.38
$ws.Range("L45").Value = 1983.5
$ws.Range("M45").Value = -262747.38
$ws.Range("N45").Value = -2737.5
$ws.Range("H61").Value = 503825.34
$ws.Range("I61").Value = 3788.6667
$ws.Range("J61").Value = 912946.25
$ws.Range("K61").Value = 3788.6667
$ws.Range("L61").Value = 912946.25
$ws.Range("M61").Value = -3576.6667
$ws.Range("N61").Value = -913370.25
$ws.Range("H74").Value = 16130698
$ws.Range("I74").Value = 1347.3529
$ws.Range("J74").Value = 35716336
$ws.Range("K74").Value = 1347.3529
$ws.Range("L74").Value = 35716336
$ws.Range("M74").Value = -473.3529000000001
$ws.Range("N74").Value = -35718084
$ws.Range("H77").Value = 16130698
$ws.Range("I77").Value = 1347.3529
$ws.Range("J77").Value = 35716336
$ws.Range("K77").Value = 6736.7645
$ws.Range("L77").Value = 178581680
$ws.Range("M77").Value = -2368.7645
$ws.Range("N77").Value = -178590416
$ws.Range("H136").Value = 503825.34
$ws.Range("I136").Value = 3788.6667
$ws.Range("J136").Value = 912946.25
$ws.Range("K136").Value = 11366.0001
$ws.Range("L136").Value = 2738838.75
$ws.Range("M136").Value = -8816.000100000001
$ws.Range("N136").Value = -2743938.75

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H31").Value = 11369546
$ws.Range("I31").Value = 1966
$ws.Range("K31").Value = 1966
$ws.Range("M31").Value = -1671
$ws.Range("H34").Value = 11369546
$ws.Range("I34").Value = 1966
$ws.Range("K34").Value = 1966
$ws.Range("M34").Value = -1764
$ws.Range("H58").Value = 8835311
$ws.Range("I58").Value = 13890992
$ws.Range("J58").Value = 1251789.4
$ws.Range("K58").Value = 13890992
$ws.Range("L58").Value = 1251789.4
$ws.Range("M58").Value = -13890789
$ws.Range("N58").Value = -1252195.4
$ws.Range("H132").Value = 5884018
$ws.Range("I132").Value = 8334737
$ws.Range("J132").Value = 2293.1
$ws.Range("K132").Value = 25004211
$ws.Range("L132").Value = 6879.299999999999
$ws.Range("M132").Value = -25001681
$ws.Range("N132").Value = -11939.3
$ws.Range("H134").Value = 12265167
$ws.Range("I134").Value = 15877049
$ws.Range("J134").Value = 1429519.9
$ws.Range("K134").Value = 47631147
$ws.Range("L134").Value = 4288559.699999999
$ws.Range("M134").Value = -47628612
$ws.Range("N134").Value = -4293629.699999999
$ws.Range("H136").Value = 8835311
$ws.Range("I136").Value = 13890992
$ws.Range("J136").Value = 1251789.4
$ws.Range("K136").Value = 41672976
$ws.Range("L136").Value = 3755368.2
$ws.Range("M136").Value = -41670426
$ws.Range("N136").Value = -3760468.2

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2115491.5
$ws.Range("J5").Value = 5684375
$ws.Range("L5").Value = 17053125
$ws.Range("N5").Value = -17053349
$ws.Range("H131").Value = 2084342.4
$ws.Range("I131").Value = 6667350.5
$ws.Range("J131").Value = 1156.7878
$ws.Range("K131").Value = 20002051.5
$ws.Range("L131").Value = 3470.3634
$ws.Range("M131").Value = -19997011.5
$ws.Range("N131").Value = -13550.3634
$ws.Range("H135").Value = 2115491.5
$ws.Range("J135").Value = 5684375
$ws.Range("L135").Value = 51159375
$ws.Range("N135").Value = -51164445

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 20000
$ws.Range("J53").Value = 20000
$ws.Range("L53").Value = 20000
$ws.Range("N53").Value = -21262
$ws.Range("H102").Value = 2741.2354
$ws.Range("I102").Value = 2174
$ws.Range("K102").Value = 2174
$ws.Range("M102").Value = -552
$ws.Range("H132").Value = 6176456.5
$ws.Range("I132").Value = 10417682
$ws.Range("J132").Value = 7400.909
$ws.Range("K132").Value = 31253046
$ws.Range("L132").Value = 22202.727
$ws.Range("M132").Value = -31250516
$ws.Range("N132").Value = -27262.727

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3140.0227
$ws.Range("I132").Value = 2901.7632
$ws.Range("J132").Value = 4649
$ws.Range("K132").Value = 8705.2896
$ws.Range("L132").Value = 13947
$ws.Range("M132").Value = -6175.2896
$ws.Range("N132").Value = -19007
$ws.Range("H134").Value = 51163.355
$ws.Range("J134").Value = 51163.355
$ws.Range("L134").Value = 51163.355
$ws.Range("N134").Value = -61303.355

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 159478.83
$ws.Range("I62").Value = 190374
$ws.Range("J62").Value = 5003
$ws.Range("K62").Value = 190374
$ws.Range("L62").Value = 5003
$ws.Range("M62").Value = -189750
$ws.Range("N62").Value = -6251
$ws.Range("H65").Value = 159478.83
$ws.Range("I65").Value = 190374
$ws.Range("J65").Value = 5003
$ws.Range("K65").Value = 951870
$ws.Range("L65").Value = 25015
$ws.Range("M65").Value = -948750
$ws.Range("N65").Value = -31255
$ws.Range("H98").Value = 43000
$ws.Range("J98").Value = 43000
$ws.Range("L98").Value = 43000
$ws.Range("N98").Value = -48990
$ws.Range("H108").Value = 42875.332
$ws.Range("J108").Value = 42875.332
$ws.Range("L108").Value = 42875.332
$ws.Range("N108").Value = -50555.332
